$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.125.60'
$ws.Range("E2").Value = '  -3.29%  '
$ws.Range("D3").Value = '1.849.03'
$ws.Range("E3").Value = '  -2.25%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7014'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -5.07%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '237.99'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.90%  '
$ws.Range("E7").Value = '  +0.13%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3030'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -4.53%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07471'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.66%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.30'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.55%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08120'
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = '1.848.00'
$ws.Range("E12").Value = '  -5.80%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.7233'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.88%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.207'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.46%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '88.59'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.86%  '
$ws.Range("D16").Value = '29.202.38'
$ws.Range("E16").Value = '  -2.91%  '
$ws.Range("E17").Value = '  -6.79%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '236.58'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -5.51%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.04'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.38%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007629'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.30%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.000'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.26%  '
$ws.Range("D22").Value = '2.108.54'
$ws.Range("E22").Value = '  +0.27%  '
$ws.Range("E23").Value = '  +0.17%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.535'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -5.03%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.961'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.79%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '160.92'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.22%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1447'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -8.62%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.00'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.08%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.957'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -5.23%  '
$ws.Range("E30").Value = '  -6.02%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.518'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.50%  '
$ws.Range("E32").Value = '  -3.12%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.953'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.97%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05130'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.41%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.182'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -6.00%  '
$ws.Range("B36").Value = 'Frax'
$ws.Range("C36").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.041'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.63%  '
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.6968'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -10.17%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.658'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.64%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01861'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.32%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.677'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.10%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9419'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +7.80%  '
$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").Value = '1.079.07'
$ws.Range("E42").Value = '  -1.72%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.944'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.21%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4268'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -6.55%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '69.55'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.50%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.000'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.09%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '101.83'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.34%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.735'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -7.03%  '
$ws.Range("D49").Value = '1.994.83'
$ws.Range("E49").Value = '  -2.75%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.157'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.97%  '
$ws.Range("B51").Value = 'Aptos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.007'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -7.78%  '
